$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = -0.0002
$ws.Range("E8").Value = 592.7477
$ws.Range("F8").Value = 592.7475
$ws.Range("D9").Value = -0.0001
$ws.Range("E9").Value = 582.0274
$ws.Range("F9").Value = 582.0273
$ws.Range("E10").Value = 477.3102
$ws.Range("F10").Value = 477.3102
$ws.Range("D11").Value = -0.00023281683255637
$ws.Range("E11").Value = 684.071315
$ws.Range("F11").Value = 684.071082183168
$ws.Range("D12").Value = -0.0000884133663471399
$ws.Range("E12").Value = 673.559924
$ws.Range("F12").Value = 673.559835586634
$ws.Range("E13").Value = 571.877265
$ws.Range("F13").Value = 571.877265
$ws.Range("D17").Value = 0.0002
$ws.Range("E17").Value = 36.5829520900008
$ws.Range("F17").Value = 36.5831520900009
$ws.Range("D18").Value = 0.0001
$ws.Range("E18").Value = 38.0333416705188
$ws.Range("F18").Value = 38.0334416705188
$ws.Range("E19").Value = 38.0754290961918
$ws.Range("F19").Value = 38.0754290961918
$ws.Range("D20").Value = 0.00023281683255637
$ws.Range("E20").Value = -54.7406629099992
$ws.Range("F20").Value = -54.7404300931668
$ws.Range("D21").Value = 0.0000884133663471399
$ws.Range("E21").Value = -53.4991823294813
$ws.Range("F21").Value = -53.4990939161149
$ws.Range("E22").Value = -56.4916359038082
$ws.Range("F22").Value = -56.4916359038082
$ws.Range("D23").Value = 0.0002
$ws.Range("E23").Value = 36.5829516950007
$ws.Range("F23").Value = 36.5831516950008
$ws.Range("D24").Value = 0.0001
$ws.Range("E24").Value = 37.9657306617447
$ws.Range("F24").Value = 37.9658306617447
$ws.Range("E25").Value = 37.5280195688705
$ws.Range("F25").Value = 37.5280195688705
$ws.Range("D26").Value = 0.00023281683255637
$ws.Range("E26").Value = -54.7406633049993
$ws.Range("F26").Value = -54.7404304881669
$ws.Range("D27").Value = 0.0000884133663471399
$ws.Range("E27").Value = -53.5667933382554
$ws.Range("F27").Value = -53.566704924889
$ws.Range("E28").Value = -57.0390454311295
$ws.Range("F28").Value = -57.0390454311295
$ws.Range("D29").Value = 0.0000328168325563703
$ws.Range("E29").Value = -91.323615
$ws.Range("F29").Value = -91.3235821831677
$ws.Range("D30").Value = -0.0000115866336528601
$ws.Range("E30").Value = -91.5325240000001
$ws.Range("F30").Value = -91.5325355866337
$ws.Range("E31").Value = -94.567065
$ws.Range("F31").Value = -94.567065
